# "wood chip as fuel"
#
# 1. Fill in the previously-blank HHV figure (20) for "wood oven dry - IPCC"
#    in row 15.
# 2. Add a new fuel row, "wood chips (EU no swiss, dry)", whose HHV/LHV/CO2
#    values are derived (by formula) from that same row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty HHV value for "wood oven dry - IPCC" (row 15).
$ws.Range("B15").Value = 20

# New row 31: wood chips (EU no swiss, dry) — references row 15's figures.
$ws.Range("A31").Value = "wood chips (EU no swiss, dry)"
$ws.Range("B31").Formula = "=B15"
$ws.Range("C31").Formula = "=C15"
$ws.Range("D31").Formula = "=D15"
